# Use inter model accuracy to ensure diverse models
# Mark the first batch of jobs (rows 2-22) in openml_100 as "Complete"
# instead of "Run", and update the active selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openml_100")

for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Value = "Complete"
}

$ws.Activate()
$ws.Range("B13").Select()

$wb.Application.Calculate()
